$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "Data" -> "Summary"
$ws.Name = "Summary"

# Remove the MSME participation breakdown row entirely
# (Micro / SMEs / MSMEs headers that lived in row 5)
$ws.Rows.Item(5).Delete()

# Add the new "title_" named style (bold + underline variant of "title")
# used by the Summary sheet template, matching the style table update.
$style = $wb.Styles.Add("title_")
$style.Font.Bold = $true
$style.Font.Underline = $true
